# Applies the "Updated cryptos list" refresh: new Price/Volume(1h) figures pulled
# from coinranking.com, plus a handful of rows whose rank (and thus row position)
# swapped with a neighbour, carrying their Coin/Link/Price/Volume along with them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. Price-column values that would otherwise be
# auto-coerced to a number by Excel (dropping e.g. a trailing '.00') are
# prefixed with a literal apostrophe so they are stored as text, matching
# the original sheet where the whole Price column holds text values.
$updates = [ordered]@{
    "D2" = "67.732.76"
    "E2" = "  +7.04%  "
    "D3" = "3.518.98"
    "E3" = "  +8.78%  "
    "E4" = "  +0.12%  "
    "D5" = "'193.08"
    "E5" = "  +9.76%  "
    "D6" = "'551.10"
    "E6" = "  +5.99%  "
    "D7" = "3.513.52"
    "E7" = "  +8.81%  "
    "D8" = "'0.608"
    "E8" = "  +2.89%  "
    "E9" = "  -0.18%  "
    "D10" = "'0.640"
    "E10" = "  +6.92%  "
    "D11" = "'57.25"
    "E11" = "  +5.74%  "
    "E12" = "  +15.78%  "
    "D13" = "'0.0000275"
    "E13" = "  +8.73%  "
    "D14" = "'9.48"
    "E14" = "  +6.08%  "
    "D15" = "4.083.70"
    "E15" = "  +9.35%  "
    "D16" = "3.521.70"
    "E16" = "  +9.13%  "
    "D17" = "67.711.79"
    "E17" = "  +7.37%  "
    "E18" = "  +4.95%  "
    "D19" = "'18.35"
    "E19" = "  +6.72%  "
    "D20" = "'11.89"
    "E20" = "  +8.53%  "
    "D21" = "'1.00"
    "E21" = "  +5.37%  "
    "D22" = "'409.91"
    "E22" = "  +11.83%  "
    "E23" = "  +6.50%  "
    "D24" = "'84.82"
    "E24" = "  +6.02%  "
    "B25" = "Toncoin"
    "C25" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    "D25" = "'4.23"
    "E25" = "  +10.47%  "
    "B26" = "RenderToken"
    "C26" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D26" = "'11.63"
    "E26" = "  +2.89%  "
    "D27" = "'2.93"
    "E27" = "  +11.86%  "
    "E28" = "  -1.28%  "
    "E29" = "  +6.17%  "
    "D30" = "'8.71"
    "E30" = "  +6.77%  "
    "D31" = "'30.68"
    "E31" = "  +8.49%  "
    "D32" = "'690.03"
    "E32" = "  +8.06%  "
    "D33" = "'6.88"
    "E33" = "  +8.02%  "
    "D34" = "'11.79"
    "E34" = "  +6.36%  "
    "D35" = "'0.112"
    "E35" = "  +7.75%  "
    "D36" = "'60.87"
    "E36" = "  +5.93%  "
    "B37" = "InjectiveProtocol"
    "C37" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D37" = "'39.36"
    "E37" = "  +7.52%  "
    "B38" = "PEPE"
    "C38" = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
    "D38" = "0.0₃0833"
    "E38" = "  +20.05%  "
    "D39" = "'0.401"
    "E39" = "  +6.37%  "
    "D40" = "'1.00"
    "E40" = "  -0.09%  "
    "D41" = "'3.41"
    "E41" = "  +25.64%  "
    "E42" = "  +10.20%  "
    "B43" = "Fetch.AI"
    "C43" = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    "D43" = "'2.74"
    "E43" = "  +9.62%  "
    "B44" = "FirstDigitalUSD"
    "C44" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "D44" = "'1.00"
    "E44" = "  +0.61%  "
    "B45" = "ThetaToken"
    "C45" = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
    "D45" = "'3.01"
    "E45" = "  +15.13%  "
    "D46" = "3.039.96"
    "E46" = "  +5.11%  "
    "B47" = "ApeXProtocol"
    "C47" = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
    "D47" = "'3.42"
    "E47" = "  +15.48%  "
    "B48" = "VeChain"
    "C48" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D48" = "'0.0426"
    "E48" = "  +9.51%  "
    "B49" = "THORChain"
    "C49" = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
    "D49" = "'9.21"
    "E49" = "  +21.81%  "
    "B50" = "WEMIXToken"
    "C50" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D50" = "'2.73"
    "E50" = "  +2.12%  "
    "E51" = "  +6.30%  "
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
